# Update the Inha -> Acvr2a ligand-receptor table with refreshed TPM-derived
# values and add the "MuSCs" sending-cluster rows (rows 5-7), matching the
# newly computed NATMI output for this sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column order: A..T =
# Sending cluster, Ligand symbol, Receptor symbol, Target cluster,
# Ligand-expressing cells, Ligand detection rate, Ligand average expression value,
# Ligand total expression value, Ligand derived specificity of average expression value,
# Ligand derived specificity of total expression value, Receptor-expressing cells,
# Receptor detection rate, Receptor average expression value, Receptor total expression value,
# Receptor derived specificity of average expression value, Receptor derived specificity of total expression value,
# Edge average expression weight, Edge total expression weight,
# Edge average expression derived specificity, Edge total expression derived specificity
$rowData = @(
    @{r=2; vals=@("FAPs", "Inha", "Acvr2a", "ECs", 1, 0.3333333333333333, 0.1623533333333333, 0.48706, 0.2351552968361266, 0.2351552968361266, 3, 1, 14.321881, 42.965643, 0.2949569176783066, 0.2949569176783066, 2.325205119953333, 20.92684607958, 0.06936068153051113, 0.06936068153051114)},
    @{r=3; vals=@("FAPs", "Inha", "Acvr2a", "FAPs", 1, 0.3333333333333333, 0.1623533333333333, 0.48706, 0.2351552968361266, 0.2351552968361266, 3, 1, 27.084169, 81.25250700000001, 0.557794259435499, 0.557794259435499, 4.397205117713334, 39.57484605942, 0.1311682746510422, 0.1311682746510422)},
    @{r=4; vals=@("FAPs", "Inha", "Acvr2a", "MuSCs", 1, 0.3333333333333333, 0.1623533333333333, 0.48706, 0.2351552968361266, 0.2351552968361266, 3, 1, 7.149790333333333, 21.449371, 0.1472488228861944, 0.1472488228861943, 1.160792293251111, 10.44713063926, 0.03462634065457326, 0.03462634065457326)},
    @{r=5; vals=@("MuSCs", "Inha", "Acvr2a", "ECs", 2, 0.6666666666666666, 0.5280556666666667, 1.584167, 0.7648447031638734, 0.7648447031638734, 3, 1, 14.321881, 42.965643, 0.2949569176783066, 0.2949569176783066, 7.562750419375667, 68.064753774381, 0.2255962361477954, 0.2255962361477954)},
    @{r=6; vals=@("MuSCs", "Inha", "Acvr2a", "FAPs", 2, 0.6666666666666666, 0.5280556666666667, 1.584167, 0.7648447031638734, 0.7648447031638734, 3, 1, 27.084169, 81.25250700000001, 0.557794259435499, 0.557794259435499, 14.30194891740767, 128.717540256669, 0.4266259847844568, 0.4266259847844568)},
    @{r=7; vals=@("MuSCs", "Inha", "Acvr2a", "MuSCs", 2, 0.6666666666666666, 0.5280556666666667, 1.584167, 0.7648447031638734, 0.7648447031638734, 3, 1, 7.149790333333333, 21.449371, 0.1472488228861944, 0.1472488228861943, 3.775487300995222, 33.979385708957, 0.1126224822316211, 0.1126224822316211)}
)

foreach ($rowDef in $rowData) {
    $r = $rowDef.r
    $vals = $rowDef.vals
    for ($i = 0; $i -lt $vals.Count; $i++) {
        $ws.Cells.Item($r, $i + 1).Value = $vals[$i]
    }
}
